$d = $word.ActiveDocument

# --- Step 1: merge the "{%p if" + " " run pairs into a single run
# "{%p if " (xml:space="preserve"). There are 3 such occurrences (one each
# highlighted yellow / cyan / magenta) scattered through the document, so
# repeat the find/replace until no more matches remain. Word's replace
# collapses the two source runs into one run (keeping the first run's
# formatting, highlight included - that gets cleared in step 2 below).
$merge = $d.Content
$merge.Find.ClearFormatting()
$merge.Find.Replacement.ClearFormatting()
$merge.Find.Text = "{%p if "
$merge.Find.Forward = $true
$merge.Find.Wrap = 0
$guard = 0
while ($merge.Find.Execute($null, $false, $false, $false, $false, $false, $true, 0, $false, "{%p if ", 2)) {
    $merge.Collapse(0)
    $guard += 1
    if ($guard -gt 50) { break }
}

# --- Step 2: strip the highlight ("clear") from every run that still has
# one - the newly merged "{%p if " runs above plus "elif", the remaining
# "{%p " runs and the "endif %}" runs. Loop using the Highlight search
# criterion (with no search text) so it keeps finding highlighted runs
# until none are left, then run one extra confirmation pass.
for ($pass = 0; $pass -lt 2; $pass++) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = ""
    $rng.Find.Highlight = $true
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 0
    $guard2 = 0
    while ($rng.Find.Execute()) {
        $rng.HighlightColorIndex = 0
        $rng.Collapse(0)
        $guard2 += 1
        if ($guard2 -gt 50) { break }
    }
}
